$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced back to
# Text format first, since these columns store formatted/text price data
# (e.g. "208.94") rather than numeric values.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.467.42'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '1.570.77'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '208.94'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D8').Value = '22.22'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').Value = '0.0593'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.793.11'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = '1.572.25'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = '63.78'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '27.478.35'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '214.25'
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').Value = '0.0₃0692'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '7.30'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '9.57'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').Value = '152.30'
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = '6.72'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -1.65%  '
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').Value = '1.382.94'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').Value = '  -0.83%  '
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').Value = '0.542'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').Value = '0.982'
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('D43').Value = '1.80'
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('D44').Value = '64.24'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = '1.705.70'
$ws.Range('D48').Value = '85.51'
$ws.Range('E48').Value = '  -3.37%  '
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').Value = '0.0498'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').Value = '0.0958'
$ws.Range('E51').Value = '  -1.46%  '
